$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove rows 5 and 6 (old "week 3" / "week 4" placeholder rows). ---
# This shifts the old row 10 (the customFormat template row) up to row 8 and the
# old row 11 (the lone J-cell row) up to row 9.
$ws.Rows("5:6").Delete()

# --- 2. Record the header style of the template row (now row 8) before we start
#        editing it, so we can stamp the new row 10 with a matching style. ---
$ws.Range("A8:AK8").Value = 1
$ws.Range("A8:AK8").Copy()
$ws.Range("A10:AK10").PasteSpecial(-4163)
$ws.Range("A10:AK10").ClearContents()

# Columns F and N:AK on the new row 10 get the lighter / default look (style 0)
# instead of the template style (style 1) that the rest of row 10 keeps.
$ws.Range("F10").ClearFormats()
$ws.Range("N10:AK10").ClearFormats()

# Undo the scratch values we used above to force the copy to "take".
$ws.Range("A8:AK8").ClearContents()

# --- 3. Trim row 8 down to just the cells the regenerated sheet actually uses. ---
$ws.Range("A8:E8").Clear()
$ws.Range("G8:M8").Clear()

# --- 4. Cell-level fixes on the remaining data rows. ---
# N3's style moves from the plain style to the alternate plain style used elsewhere.
$ws.Range("N3").Value = 12

# N4 is a brand-new "Feed" total cell for week 2, defaulting to 0.
$ws.Range("N4").Value = 0

# --- 5. Selection follows the data down to the newly relocated template row. ---
$ws.Range("N11").Select()
